$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 2000
$ws.Range("I21").Value = 2000
$ws.Range("K21").Value = 2000
$ws.Range("M21").Value = -1532
$ws.Range("H23").Value = 2000
$ws.Range("I23").Value = 2000
$ws.Range("K23").Value = 2000
$ws.Range("M23").Value = -1766
$ws.Range("H55").Value = 83
$ws.Range("J55").Value = 85
$ws.Range("L55").Value = 85
$ws.Range("N55").Value = -513
$ws.Range("H87").Value = 34999.25
$ws.Range("J87").Value = 34999.25
$ws.Range("L87").Value = 34999.25
$ws.Range("N87").Value = -37495.25
$ws.Range("H90").Value = 34999.25
$ws.Range("J90").Value = 34999.25
$ws.Range("L90").Value = 104997.75
$ws.Range("N90").Value = -117477.75
$ws.Range("H112").Value = 2302.889
$ws.Range("I112").Value = 1649.5
$ws.Range("K112").Value = 4948.5
$ws.Range("M112").Value = -3840.5
$ws.Range("H138").Value = 4244.5
$ws.Range("I138").Value = 2194.6667
$ws.Range("J138").Value = 4927.778
$ws.Range("K138").Value = 6584.000100000001
$ws.Range("L138").Value = 14783.334
$ws.Range("M138").Value = -1444.000100000001
$ws.Range("N138").Value = -25063.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H86").Value = 314314
$ws.Range("J86").Value = 314314
$ws.Range("L86").Value = 314314
$ws.Range("N86").Value = -316686
$ws.Range("H88").Value = 2327.4285
$ws.Range("I88").Value = 1364.6666
$ws.Range("J88").Value = 3049.5
$ws.Range("K88").Value = 1364.6666
$ws.Range("L88").Value = 3049.5
$ws.Range("M88").Value = -958.6666
$ws.Range("N88").Value = -3861.5
$ws.Range("H89").Value = 314314
$ws.Range("J89").Value = 314314
$ws.Range("L89").Value = 942942
$ws.Range("N89").Value = -954798
$ws.Range("H91").Value = 2327.4285
$ws.Range("I91").Value = 1364.6666
$ws.Range("J91").Value = 3049.5
$ws.Range("K91").Value = 1364.6666
$ws.Range("L91").Value = 3049.5
$ws.Range("M91").Value = 39.33339999999998
$ws.Range("N91").Value = -5857.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 253.66667
$ws.Range("I11").Value = 333.25
$ws.Range("J11").Value = 94.5
$ws.Range("K11").Value = 333.25
$ws.Range("L11").Value = 94.5
$ws.Range("M11").Value = -193.25
$ws.Range("N11").Value = -374.5
$ws.Range("H35").Value = 14998.333
$ws.Range("J35").Value = 14998.333
$ws.Range("L35").Value = 14998.333
$ws.Range("N35").Value = -15618.333
$ws.Range("H94").Value = 2194
$ws.Range("I94").Value = 1771.4736
$ws.Range("K94").Value = 1771.4736
$ws.Range("M94").Value = -1320.4736
$ws.Range("H99").Value = 4450
$ws.Range("I99").Value = 4450
$ws.Range("K99").Value = 4450
$ws.Range("M99").Value = -2952
$ws.Range("H105").Value = 2794.75
$ws.Range("I105").Value = 1926.3334
$ws.Range("K105").Value = 1926.3334
$ws.Range("M105").Value = -179.3334
$ws.Range("H134").Value = 5300
$ws.Range("I134").Value = 5052.75
$ws.Range("K134").Value = 15158.25
$ws.Range("M134").Value = -12623.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2750
$ws.Range("I6").Value = 2333.3333
$ws.Range("K6").Value = 2333.3333
$ws.Range("M6").Value = -2220.3333
$ws.Range("H10").Value = 166.66667
$ws.Range("I10").Value = 166.66667
$ws.Range("K10").Value = 166.66667
$ws.Range("M10").Value = -27.66667000000001
$ws.Range("H31").Value = 4203.1
$ws.Range("I31").Value = 3799.5
$ws.Range("J31").Value = 4472.1665
$ws.Range("K31").Value = 3799.5
$ws.Range("L31").Value = 4472.1665
$ws.Range("M31").Value = -3504.5
$ws.Range("N31").Value = -5062.1665
$ws.Range("H34").Value = 4203.1
$ws.Range("I34").Value = 3799.5
$ws.Range("J34").Value = 4472.1665
$ws.Range("K34").Value = 3799.5
$ws.Range("L34").Value = 4472.1665
$ws.Range("M34").Value = -3597.5
$ws.Range("N34").Value = -4876.1665
$ws.Range("H58").Value = 9054.75
$ws.Range("I58").Value = 7106.8
$ws.Range("J58").Value = 12301.333
$ws.Range("K58").Value = 7106.8
$ws.Range("L58").Value = 12301.333
$ws.Range("M58").Value = -6903.8
$ws.Range("N58").Value = -12707.333
$ws.Range("H59").Value = 30802.062
$ws.Range("I59").Value = 23809.666
$ws.Range("K59").Value = 23809.666
$ws.Range("M59").Value = -22664.666
$ws.Range("H60").Value = 17187
$ws.Range("I60").Value = 7422
$ws.Range("K60").Value = 7422
$ws.Range("M60").Value = -6911
$ws.Range("H107").Value = 614.4
$ws.Range("I107").Value = 532
$ws.Range("K107").Value = 532
$ws.Range("M107").Value = 1388
$ws.Range("H136").Value = 9054.75
$ws.Range("I136").Value = 7106.8
$ws.Range("J136").Value = 12301.333
$ws.Range("K136").Value = 21320.4
$ws.Range("L136").Value = 36903.999
$ws.Range("M136").Value = -18770.4
$ws.Range("N136").Value = -42003.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 100
$ws.Range("I24").Value = 100
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 300
$ws.Range("L24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -70
$ws.Range("H57").Value = 5226
$ws.Range("I57").Value = 3543.3333
$ws.Range("J57").Value = 7750
$ws.Range("K57").Value = 10629.9999
$ws.Range("L57").Value = 23250
$ws.Range("M57").Value = -10070.9999
$ws.Range("N57").Value = -24368
$ws.Range("H81").Value = 6130.1333
$ws.Range("J81").Value = 6815.636
$ws.Range("L81").Value = 20446.908
$ws.Range("N81").Value = -22692.908
$ws.Range("H84").Value = 6130.1333
$ws.Range("J84").Value = 6815.636
$ws.Range("L84").Value = 61340.724
$ws.Range("N84").Value = -72572.724
$ws.Range("H113").Value = 1404
$ws.Range("J113").Value = 1404
$ws.Range("L113").Value = 4212
$ws.Range("N113").Value = -8552
$ws.Range("H128").Value = 342000
$ws.Range("I128").Value = 342000
$ws.Range("K128").Value = 1026000
$ws.Range("M128").Value = -1021020
$ws.Range("H140").Value = 3126.1667
$ws.Range("I140").Value = 2702.6667
$ws.Range("K140").Value = 8108.000100000001
$ws.Range("M140").Value = -2928.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 2518.7
$ws.Range("I9").Value = 623.375
$ws.Range("J9").Value = 10100
$ws.Range("K9").Value = 623.375
$ws.Range("L9").Value = 10100
$ws.Range("M9").Value = -453.375
$ws.Range("N9").Value = -10440
$ws.Range("H57").Value = 25000
$ws.Range("J57").Value = 25000
$ws.Range("L57").Value = 25000
$ws.Range("N57").Value = -26640
$ws.Range("H62").Value = 34000
$ws.Range("J62").Value = 34000
$ws.Range("L62").Value = 34000
$ws.Range("N62").Value = -35372
$ws.Range("H65").Value = 34000
$ws.Range("J65").Value = 34000
$ws.Range("L65").Value = 102000
$ws.Range("N65").Value = -108864
$ws.Range("H80").Value = 5833.1665
$ws.Range("I80").Value = 6249.75
$ws.Range("K80").Value = 6249.75
$ws.Range("M80").Value = -5251.75
$ws.Range("H83").Value = 5833.1665
$ws.Range("I83").Value = 6249.75
$ws.Range("K83").Value = 31248.75
$ws.Range("M83").Value = -26256.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1689.1
$ws.Range("I7").Value = 1689.1
$ws.Range("K7").Value = 1689.1
$ws.Range("M7").Value = -1577.1
$ws.Range("H9").Value = 1374.25
$ws.Range("J9").Value = 3999
$ws.Range("L9").Value = 3999
$ws.Range("N9").Value = -4447
$ws.Range("H46").Value = 2209.875
$ws.Range("I46").Value = 2209.875
$ws.Range("K46").Value = 2209.875
$ws.Range("M46").Value = -2021.875
$ws.Range("H93").Value = 2157.6667
$ws.Range("I93").Value = 1668.75
$ws.Range("J93").Value = 6069
$ws.Range("K93").Value = 1668.75
$ws.Range("L93").Value = 6069
$ws.Range("M93").Value = -420.75
$ws.Range("N93").Value = -8565
$ws.Range("H126").Value = 1689.1
$ws.Range("I126").Value = 1689.1
$ws.Range("K126").Value = 5067.299999999999
$ws.Range("M126").Value = -2597.299999999999
$ws.Range("H127").Value = 80000
$ws.Range("J127").Value = 80000
$ws.Range("L127").Value = 80000
$ws.Range("N127").Value = -89920

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1041.6666
$ws.Range("I81").Value = 937.5
$ws.Range("K81").Value = 1875
$ws.Range("M81").Value = -814
$ws.Range("H84").Value = 1041.6666
$ws.Range("I84").Value = 937.5
$ws.Range("K84").Value = 9375
$ws.Range("M84").Value = -4071
$ws.Range("H107").Value = 407.5
$ws.Range("I107").Value = 411.25
$ws.Range("K107").Value = 1233.75
$ws.Range("M107").Value = 686.25
$ws.Range("H122").Value = 3859.8572
$ws.Range("I122").Value = 1029.75
$ws.Range("J122").Value = 7633.3335
$ws.Range("K122").Value = 3089.25
$ws.Range("L122").Value = 22900.0005
$ws.Range("M122").Value = -639.25
$ws.Range("N122").Value = -27800.0005
